$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "S" marker in the newly-checked cells (Reuniao 24/09 and 25/09 tracking columns)
$ws.Range("J6").Value = "S"
$ws.Range("I7").Value = "S"
$ws.Range("J7").Value = "S"
$ws.Range("K7").Value = "S"
$ws.Range("L7").Value = "S"
$ws.Range("L8").Value = "S"
$ws.Range("I9").Value = "S"
$ws.Range("L10").Value = "S"
$ws.Range("I11").Value = "S"
$ws.Range("J14").Value = "S"
$ws.Range("K14").Value = "S"
$ws.Range("L14").Value = "S"

# Update selection to match the new active cell
$ws.Range("L9").Select()
